# Apply cryptos-list price/volume refresh (GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.132.42"

# Row 3
$ws.Range("D3").Value = "3.116.09"
$ws.Range("E3").Value = "  +4.00%  "

# Row 4
$ws.Range("E4").Value = "  -0.15%  "

# Row 5
$ws.Range("D5").Value = "'585.62"
$ws.Range("E5").Value = "  +4.00%  "

# Row 6
$ws.Range("D6").Value = "'144.76"
$ws.Range("E6").Value = "  +4.18%  "

# Row 7
$ws.Range("E7").Value = "  -0.11%  "

# Row 8
$ws.Range("D8").Value = "3.109.36"
$ws.Range("E8").Value = "  +4.40%  "

# Row 9
$ws.Range("E9").Value = "  +1.80%  "

# Row 10
$ws.Range("E10").Value = "  +13.10%  "

# Row 11
$ws.Range("D11").Value = "'5.80"
$ws.Range("E11").Value = "  +10.00%  "

# Row 12
$ws.Range("E12").Value = "  +3.10%  "

# Row 13
$ws.Range("D13").Value = "'0.0000248"
$ws.Range("E13").Value = "  +8.18%  "

# Row 14
$ws.Range("D14").Value = "'35.54"
$ws.Range("E14").Value = "  +5.25%  "

# Row 15
$ws.Range("E15").Value = "  +0.35%  "

# Row 16
$ws.Range("D16").Value = "3.630.06"
$ws.Range("E16").Value = "  +4.03%  "

# Row 17
$ws.Range("D17").Value = "'7.17"
$ws.Range("E17").Value = "  -0.07%  "

# Row 18
$ws.Range("D18").Value = "63.032.82"
$ws.Range("E18").Value = "  +5.95%  "

# Row 19
$ws.Range("D19").Value = "3.111.52"
$ws.Range("E19").Value = "  +4.00%  "

# Row 20
$ws.Range("D20").Value = "'467.77"
$ws.Range("E20").Value = "  +7.59%  "

# Row 21
$ws.Range("D21").Value = "'14.09"
$ws.Range("E21").Value = "  +3.94%  "

# Row 22
$ws.Range("D22").Value = "'0.726"
$ws.Range("E22").Value = "  +1.24%  "

# Row 23
$ws.Range("E23").Value = "  +7.16%  "

# Row 24
$ws.Range("D24").Value = "'13.29"
$ws.Range("E24").Value = "  -1.09%  "

# Row 25
$ws.Range("E25").Value = "  +2.20%  "

# Row 26
$ws.Range("E26").Value = "  +0.06%  "

# Row 27
$ws.Range("D27").Value = "'8.39"
$ws.Range("E27").Value = "  +8.08%  "

# Row 28
$ws.Range("E28").Value = "  +0.59%  "

# Row 29
$ws.Range("E29").Value = "  +5.45%  "

# Row 30
$ws.Range("E30").Value = "  -0.18%  "

# Row 31
$ws.Range("D31").Value = "'6.83"
$ws.Range("E31").Value = "  +9.56%  "

# Row 32
$ws.Range("D32").Value = "'26.94"
$ws.Range("E32").Value = "  +4.67%  "

# Row 33
$ws.Range("D33").Value = "'0.110"
$ws.Range("E33").Value = "  +3.49%  "

# Row 34
$ws.Range("D34").Value = "0.0₃0861"
$ws.Range("E34").Value = "  +10.59%  "

# Row 35
$ws.Range("E35").Value = "  +15.52%  "

# Row 36
$ws.Range("E36").Value = "  +5.09%  "

# Row 37
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").Value = "'3.31"
$ws.Range("E37").Value = "  +19.25%  "

# Row 38
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").Value = "'6.02"
$ws.Range("E38").Value = "  +2.48%  "

# Row 39
$ws.Range("D39").Value = "'50.89"
$ws.Range("E39").Value = "  +4.14%  "

# Row 40
$ws.Range("D40").Value = "'430.55"
$ws.Range("E40").Value = "  +7.46%  "

# Row 41
$ws.Range("D41").Value = "'8.72"
$ws.Range("E41").Value = "  +1.50%  "

# Row 42
$ws.Range("D42").Value = "2.927.48"
$ws.Range("E42").Value = "  +6.07%  "

# Row 43
$ws.Range("E43").Value = "  +4.41%  "

# Row 44
$ws.Range("E44").Value = "  +11.10%  "

# Row 45
$ws.Range("E45").Value = "  +5.57%  "

# Row 46
$ws.Range("E46").Value = "  +7.99%  "

# Row 47
$ws.Range("E47").Value = "  +2.77%  "

# Row 49
$ws.Range("D49").Value = "'123.55"
$ws.Range("E49").Value = "  +0.52%  "

# Row 50
$ws.Range("E50").Value = "  +0.87%  "

# Row 51
$ws.Range("D51").Value = "'24.53"
$ws.Range("E51").Value = "  +4.30%  "

